$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.5814128013667214
$ws.Range("D2").Value = 0.5668722678622609

$ws.Range("C3").Value = 1.376583167138703
$ws.Range("D3").Value = 0.1824907262166702

$ws.Range("C4").Value = 1.691001067455943
$ws.Range("D4").Value = 0.1049565830253625

$ws.Range("C5").Value = 0.3765569824416972
$ws.Range("D5").Value = 0.7101101931556548

$ws.Range("C6").Value = 0.6288778141945918
$ws.Range("D6").Value = 0.5359047978319764

$ws.Range("C7").Value = 0.8939663943343991
$ws.Range("D7").Value = 0.3810158041978626

$ws.Range("C8").Value = -0.1891148604798456
$ws.Range("D8").Value = 0.8517366148365066

$ws.Range("C9").Value = 0.2696952561483115
$ws.Range("D9").Value = 0.7899082369468893

$ws.Range("C10").Value = -1.083825304400224
$ws.Range("D10").Value = 0.2901778663374535

$ws.Range("C11").Value = -1.196327754742788
$ws.Range("D11").Value = 0.2443031152981359
